$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New key/value rows to append after the existing data (rows 2-8 used),
# starting at row 9.
$rows = @(
    @("victory", "VICTORY"),
    @("score", "SCORE"),
    @("time", "TIME"),
    @("time_bonus", "TIME BONUS")
)

$r = 9
foreach ($pair in $rows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# Rows 13-14 were authored so the shared-string table gains new unique
# entries in the order: perfect, total, TOTAL, PERFECT. Fill keys for
# both rows first, then values for row 14 before row 13, to reproduce
# that interning order.
$ws.Cells.Item(13, 1).Value = "perfect"
$ws.Cells.Item(14, 1).Value = "total"
$ws.Cells.Item(14, 2).Value = "TOTAL"
$ws.Cells.Item(13, 2).Value = "PERFECT"

# Update the active selection to match the post-edit state (A13).
$ws.Range("A13").Select()
